$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = 46027
$ws.Range("B2").Value = 14306.1587118229
$ws.Range("C2").Value = 13112.2375494424
$ws.Range("D2").Value = 21158.86
$ws.Range("E2").Value = 9788.63399182388
$ws.Range("F2").Value = 72.5838142194289

$ws.Range("A3").Value = 46028
$ws.Range("B3").Value = 5967.75786836631
$ws.Range("C3").Value = 9228.9186673625
$ws.Range("E3").Value = 9473.69119614784
$ws.Range("F3").Value = 276.114577646264

$ws.Range("A4").Value = 46029
$ws.Range("B4").Value = 14362.4754236912
$ws.Range("C4").Value = 13968.0826086507
$ws.Range("E4").Value = 9893.1930333134
$ws.Range("F4").Value = 491.058985081835

$ws.Range("A5").Value = 46030
$ws.Range("B5").Value = 14412.725245648
$ws.Range("C5").Value = 14214.5795224253
$ws.Range("E5").Value = 9966.82945568614
$ws.Range("F5").Value = 504.397874087976

$ws.Range("A6").Value = 46031
$ws.Range("B6").Value = 13853.5740879566
$ws.Range("C6").Value = 12613.2955939022
$ws.Range("E6").Value = 9581.77572469658
$ws.Range("F6").Value = 421.633804941617

$ws.Range("A7").Value = 46032
$ws.Range("B7").Value = 5459.74416686733
$ws.Range("C7").Value = 8542.45579740745
$ws.Range("E7").Value = 8698.36715313551
$ws.Range("F7").Value = 215.206789605956

$ws.Range("A8").Value = 46033
$ws.Range("B8").Value = 5293.94498459496
$ws.Range("C8").Value = 8742.65655772017
$ws.Range("E8").Value = 8686.09422494168
$ws.Range("F8").Value = 223.037115944244

$ws.Range("A9").Value = 46034
$ws.Range("B9").Value = 13378.423123507
$ws.Range("C9").Value = 12486.8178951059
$ws.Range("E9").Value = 9104.8400452809
$ws.Range("F9").Value = 396.49158084945

$ws.Range("A10").Value = 46035
$ws.Range("B10").Value = 13378.423123507
$ws.Range("C10").Value = 12484.8685332908
$ws.Range("E10").Value = 9104.8400452809
$ws.Range("F10").Value = 396.410357440486

$ws.Range("A11").Value = 46036
$ws.Range("B11").Value = 13378.423123507
$ws.Range("C11").Value = 12197.246597475
$ws.Range("E11").Value = 9104.8400452809
$ws.Range("F11").Value = 384.426110114828

$ws.Range("A12").Value = 46037
$ws.Range("B12").Value = 13378.423123507
$ws.Range("C12").Value = 12255.0940623911
$ws.Range("E12").Value = 9104.8400452809
$ws.Range("F12").Value = 386.836421153

$ws.Range("A13").Value = 46038
$ws.Range("B13").Value = 13378.423123507
$ws.Range("C13").Value = 11589.2979576
$ws.Range("E13").Value = 9104.79996470691
$ws.Range("F13").Value = 359.093246762787

$ws.Range("A14").Value = 46039
$ws.Range("B14").Value = 5471.00037786234
$ws.Range("C14").Value = 7955.93762887441
$ws.Range("E14").Value = 8689.9890489968
$ws.Range("F14").Value = 190.4194449113

$ws.Range("A15").Value = 46040
$ws.Range("B15").Value = 5302.34398887746
$ws.Range("C15").Value = 8117.45551891526
$ws.Range("E15").Value = 8677.57315317719
$ws.Range("F15").Value = 196.632028003852

